# First Commit of Actual Program
# - Remove Sheet2 and Sheet3, keep only the first sheet.
# - Rename the remaining sheet ("Sheet1") to "ValidLogin".
# - Populate a small login-form style table:
#     A1: Username   B1: Password
#     A2: ADMIN      B2: manager
# - Leave the selection on A3, matching the authored workbook.

$wb = $excel.ActiveWorkbook

# Drop the extra sheets (Sheet2, Sheet3) so only one sheet remains.
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

# Rename the surviving sheet.
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "ValidLogin"

# Fill in the header row and the data row.
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "ADMIN"
$ws.Range("B2").Value = "manager"

# Match the saved selection state (A3 selected).
$ws.Range("A3").Select()
